$d = $word.ActiveDocument

# The second table holds the DATE / PROJECT NAME header rows.
# Row 1, Column 2 is the empty cell next to "DATE" that should receive the date text.
$table = $d.Tables.Item(2)
$cell = $table.Cell(1, 2)
$rng = $cell.Range

# Insert the date text into the (currently empty) cell paragraph.
$rng.Text = "22/10/22"

# Match the formatting used elsewhere on this page (Bodoni MT Black, 36pt / sz 72 half-points).
$rng.Font.Name = "Bodoni MT Black"
$rng.Font.Size = 36
$rng.Font.SizeBi = 36
